$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("geometry")

# Insert a new column before AN, shifting thickness_max_chord_ratio .. diameter_le_chord_ratio
# (old AN:AQ) one column to the right (new AO:AR).
$ws.Range("AN1").EntireColumn.Insert()

# Populate the newly inserted column with the "solidity" data.
$ws.Range("AN1").Value = "solidity"
$ws.Range("AN2").Value = "[1.42997704 1.70997375]"
